# Update the Metadata sheet:
#  - Version bumped 5.0.0 -> 6.0.0
#  - Date bumped to the new publish timestamp
#  - Publisher value filled in ("Alvearie Team")
#  - The two duplicate "Contact / No display for ContactDetail" rows are
#    replaced by a single "Jurisdiction / United States of America" row
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: row 3, column B
$ws.Range("B3").Value = "6.0.0"

# Date: row 8, column B
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher: row 9, column B (was empty)
$ws.Range("B9").Value = "Alvearie Team"

# Row 10 was "Contact" / "No display for ContactDetail" -> becomes "Jurisdiction" / "United States of America"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row -> remove it entirely,
# shifting everything below it up by one row.
$ws.Range("A11").EntireRow.Delete()
